$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row 1: new date columns AD1:AF1
$ws.Range("AD1").Value = "31/12/2023"
$ws.Range("AE1").Value = "31/03/2024"
$ws.Range("AF1").Value = "30/06/2024"
$ws.Range("AD1:AF1").Font.Bold = $true
$ws.Range("AD1:AF1").HorizontalAlignment = -4108
$ws.Range("AD1:AF1").VerticalAlignment = -4160
$ws.Range("AD1:AF1").Borders.LineStyle = 1

$ws.Range("AD2").Value = 4117466.112
$ws.Range("AE2").Value = 4590901.76
$ws.Range("AF2").Value = 4331321.856
$ws.Range("AD3").Value = 902076.992
$ws.Range("AE3").Value = 1430381.056
$ws.Range("AF3").Value = 1194092.032
$ws.Range("AD4").Value = 81279
$ws.Range("AE4").Value = 73419
$ws.Range("AF4").Value = 31843
$ws.Range("AD5").Value = 355300
$ws.Range("AE5").Value = 894598.976
$ws.Range("AF5").Value = 656443.008
$ws.Range("AD6").Value = 213712
$ws.Range("AE6").Value = 226792
$ws.Range("AF6").Value = 233656
$ws.Range("AD7").Value = 161098
$ws.Range("AE7").Value = 147696.992
$ws.Range("AF7").Value = 169488.992
$ws.Range("AD8").Value = 0
$ws.Range("AE8").Value = 0
$ws.Range("AF8").Value = 0
$ws.Range("AD9").Value = 52941
$ws.Range("AE9").Value = 53986
$ws.Range("AF9").Value = 63830
$ws.Range("AD10").Value = 0
$ws.Range("AE10").Value = 0
$ws.Range("AF10").Value = 0
$ws.Range("AD11").Value = 37747
$ws.Range("AE11").Value = 33888
$ws.Range("AF11").Value = 38831
$ws.Range("AD12").Value = 282476.992
$ws.Range("AE12").Value = 279982.016
$ws.Range("AF12").Value = 292544
$ws.Range("AD13").Value = 0
$ws.Range("AE13").Value = 1053
$ws.Range("AF13").Value = 480
$ws.Range("AD14").Value = 0
$ws.Range("AE14").Value = 0
$ws.Range("AF14").Value = 0
$ws.Range("AD15").Value = 0
$ws.Range("AE15").Value = 0
$ws.Range("AF15").Value = 0
$ws.Range("AD16").Value = 0
$ws.Range("AE16").Value = 0
$ws.Range("AF16").Value = 0
$ws.Range("AD17").Value = 0
$ws.Range("AE17").Value = 0
$ws.Range("AF17").Value = 0
$ws.Range("AD18").Value = 0
$ws.Range("AE18").Value = 0
$ws.Range("AF18").Value = 0
$ws.Range("AD19").Value = 0
$ws.Range("AE19").Value = 0
$ws.Range("AF19").Value = 0
$ws.Range("AD20").Value = 0
$ws.Range("AE20").Value = 0
$ws.Range("AF20").Value = 0
$ws.Range("AD21").Value = 0
$ws.Range("AE21").Value = 0
$ws.Range("AF21").Value = 0
$ws.Range("AD22").Value = 0
$ws.Range("AE22").Value = 0
$ws.Range("AF22").Value = 0
$ws.Range("AD23").Value = 1378694.016
$ws.Range("AE23").Value = 1378985.984
$ws.Range("AF23").Value = 1352973.952
$ws.Range("AD24").Value = 1554217.984
$ws.Range("AE24").Value = 1501553.024
$ws.Range("AF24").Value = 1491712
$ws.Range("AD25").Value = 0
$ws.Range("AE25").Value = 0
$ws.Range("AF25").Value = 0
$ws.Range("AD26").Value = 4117466.112
$ws.Range("AE26").Value = 4590901.76
$ws.Range("AF26").Value = 4331321.856
$ws.Range("AD27").Value = 1306667.008
$ws.Range("AE27").Value = 1247815.04
$ws.Range("AF27").Value = 1162528
$ws.Range("AD28").Value = 121685
$ws.Range("AE28").Value = 139320
$ws.Range("AF28").Value = 152536
$ws.Range("AD29").Value = 357539.008
$ws.Range("AE29").Value = 274508.992
$ws.Range("AF29").Value = 300300.992
$ws.Range("AD30").Value = 30685
$ws.Range("AE30").Value = 26937
$ws.Range("AF30").Value = 32242
$ws.Range("AD31").Value = 543353.024
$ws.Range("AE31").Value = 575190.016
$ws.Range("AF31").Value = 454436.992
$ws.Range("AD32").Value = 0
$ws.Range("AE32").Value = 0
$ws.Range("AF32").Value = 0
$ws.Range("AD33").Value = 0
$ws.Range("AE33").Value = 0
$ws.Range("AF33").Value = 0
$ws.Range("AD34").Value = 253404.992
$ws.Range("AE34").Value = 231859.008
$ws.Range("AF34").Value = 223012
$ws.Range("AD35").Value = 0
$ws.Range("AE35").Value = 0
$ws.Range("AF35").Value = 0
$ws.Range("AD36").Value = 0
$ws.Range("AE36").Value = 0
$ws.Range("AF36").Value = 0
$ws.Range("AD37").Value = 1417118.976
$ws.Range("AE37").Value = 2032646.016
$ws.Range("AF37").Value = 1917885.056
$ws.Range("AD38").Value = 573460.992
$ws.Range("AE38").Value = 1216188.032
$ws.Range("AF38").Value = 1114950.016
$ws.Range("AD39").Value = 0
$ws.Range("AE39").Value = 0
$ws.Range("AF39").Value = 0
$ws.Range("AD40").Value = 743681.024
$ws.Range("AE40").Value = 703104
$ws.Range("AF40").Value = 680793.984
$ws.Range("AD41").Value = 55953
$ws.Range("AE41").Value = 70500
$ws.Range("AF41").Value = 80372
$ws.Range("AD42").Value = 0
$ws.Range("AE42").Value = 0
$ws.Range("AF42").Value = 0
$ws.Range("AD43").Value = 44024
$ws.Range("AE43").Value = 42854
$ws.Range("AF43").Value = 41769
$ws.Range("AD44").Value = 0
$ws.Range("AE44").Value = 0
$ws.Range("AF44").Value = 0
$ws.Range("AD45").Value = 0
$ws.Range("AE45").Value = 0
$ws.Range("AF45").Value = 0
$ws.Range("AD46").Value = 0
$ws.Range("AE46").Value = 0
$ws.Range("AF46").Value = 0
$ws.Range("AD47").Value = 1393680
$ws.Range("AE47").Value = 1310440.96
$ws.Range("AF47").Value = 1250909.056
$ws.Range("AD48").Value = 1461068.032
$ws.Range("AE48").Value = 1461068.032
$ws.Range("AF48").Value = 1461068.032
$ws.Range("AD49").Value = 671046.976
$ws.Range("AE49").Value = 687057.984
$ws.Range("AF49").Value = 687057.984
$ws.Range("AD50").Value = 0
$ws.Range("AE50").Value = 0
$ws.Range("AF50").Value = 0
$ws.Range("AD51").Value = 0
$ws.Range("AE51").Value = 0
$ws.Range("AF51").Value = 0
$ws.Range("AD52").Value = -738419.008
$ws.Range("AE52").Value = -829187.968
$ws.Range("AF52").Value = -856691.008
$ws.Range("AD53").Value = 0
$ws.Range("AE53").Value = 0
$ws.Range("AF53").Value = 0
$ws.Range("AD54").Value = 0
$ws.Range("AE54").Value = 0
$ws.Range("AF54").Value = 0
$ws.Range("AD55").Value = -16
$ws.Range("AE55").Value = -8497
$ws.Range("AF55").Value = -40526
$ws.Range("AD56").Value = 0
$ws.Range("AE56").Value = 0
$ws.Range("AF56").Value = 0
$ws.Range("AD57:AF57").Font.Bold = $false
$ws.Range("AD58:AF58").Font.Bold = $false
$ws.Range("AD59").Value = 1074050.944
$ws.Range("AE59").Value = 1028624
$ws.Range("AF59").Value = 1108247.04
$ws.Range("AD60").Value = -317320.96
$ws.Range("AE60").Value = -369641.984
$ws.Range("AF60").Value = -387671.008
$ws.Range("AD61").Value = 756729.984
$ws.Range("AE61").Value = 658982.016
$ws.Range("AF61").Value = 720576
$ws.Range("AD62").Value = -609381.952
$ws.Range("AE62").Value = -592814.976
$ws.Range("AF62").Value = -629606.016
$ws.Range("AD63").Value = -61531
$ws.Range("AE63").Value = -99709
$ws.Range("AF63").Value = -63307
$ws.Range("AD64").Value = 0
$ws.Range("AE64").Value = 0
$ws.Range("AF64").Value = 0
$ws.Range("AD65").Value = 0
$ws.Range("AE65").Value = 0
$ws.Range("AF65").Value = 0
$ws.Range("AD66").Value = 0
$ws.Range("AE66").Value = 0
$ws.Range("AF66").Value = 0
$ws.Range("AD67").Value = 0
$ws.Range("AE67").Value = 0
$ws.Range("AF67").Value = 0
$ws.Range("AD68").Value = -39636.008
$ws.Range("AE68").Value = -42681
$ws.Range("AF68").Value = -45294
$ws.Range("AD69").Value = 10703
$ws.Range("AE69").Value = 15314
$ws.Range("AF69").Value = 16912
$ws.Range("AD70").Value = -50339
$ws.Range("AE70").Value = -57995
$ws.Range("AF70").Value = -62206
$ws.Range("AD71:AF71").Font.Bold = $false
$ws.Range("AD72:AF72").Font.Bold = $false
$ws.Range("AD73:AF73").Font.Bold = $false
$ws.Range("AD74").Value = 46181
$ws.Range("AE74").Value = -76223
$ws.Range("AF74").Value = -17631
$ws.Range("AD75").Value = 0
$ws.Range("AE75").Value = 0
$ws.Range("AF75").Value = 0
$ws.Range("AD76").Value = 13131
$ws.Range("AE76").Value = -14546
$ws.Range("AF76").Value = -9872
$ws.Range("AD77:AF77").Font.Bold = $false
$ws.Range("AD78:AF78").Font.Bold = $false
$ws.Range("AD79").Value = 0
$ws.Range("AE79").Value = 0
$ws.Range("AF79").Value = 0
$ws.Range("AD80").Value = 59312
$ws.Range("AE80").Value = -90769
$ws.Range("AF80").Value = -27503
